# edit.ps1
# Updates "想去人数" (column F) values across all 4 worksheets
# to match the data refresh captured in the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 7747
$ws.Cells.Item(4, 6).Value = 7909
$ws.Cells.Item(8, 6).Value = 6721
$ws.Cells.Item(9, 6).Value = 6721
$ws.Cells.Item(10, 6).Value = 3398
$ws.Cells.Item(12, 6).Value = 3736
$ws.Cells.Item(15, 6).Value = 46
$ws.Cells.Item(16, 6).Value = 71
$ws.Cells.Item(17, 6).Value = 77
$ws.Cells.Item(20, 6).Value = 50
$ws.Cells.Item(21, 6).Value = 325
$ws.Cells.Item(23, 6).Value = 334
$ws.Cells.Item(24, 6).Value = 3881
$ws.Cells.Item(26, 6).Value = 373
$ws.Cells.Item(27, 6).Value = 957
$ws.Cells.Item(29, 6).Value = 1501
$ws.Cells.Item(31, 6).Value = 62
$ws.Cells.Item(32, 6).Value = 2775
$ws.Cells.Item(33, 6).Value = 1897
$ws.Cells.Item(34, 6).Value = 37
$ws.Cells.Item(35, 6).Value = 52
$ws.Cells.Item(36, 6).Value = 64
$ws.Cells.Item(37, 6).Value = 57
$ws.Cells.Item(38, 6).Value = 3744
$ws.Cells.Item(39, 6).Value = 334
$ws.Cells.Item(40, 6).Value = 283
$ws.Cells.Item(41, 6).Value = 45
$ws.Cells.Item(42, 6).Value = 924
$ws.Cells.Item(43, 6).Value = 551
$ws.Cells.Item(45, 6).Value = 1445
$ws.Cells.Item(47, 6).Value = 4
$ws.Cells.Item(48, 6).Value = 561
$ws.Cells.Item(49, 6).Value = 652

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 27
$ws.Cells.Item(6, 6).Value = 416
$ws.Cells.Item(7, 6).Value = 44
$ws.Cells.Item(9, 6).Value = 103
$ws.Cells.Item(11, 6).Value = 40
$ws.Cells.Item(17, 6).Value = 121

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 139

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 13
$ws.Cells.Item(3, 6).Value = 139
$ws.Cells.Item(5, 6).Value = 27
$ws.Cells.Item(7, 6).Value = 7747
$ws.Cells.Item(9, 6).Value = 7909
$ws.Cells.Item(12, 6).Value = 6721
$ws.Cells.Item(13, 6).Value = 3398
$ws.Cells.Item(15, 6).Value = 3736
$ws.Cells.Item(18, 6).Value = 46
$ws.Cells.Item(19, 6).Value = 71
$ws.Cells.Item(20, 6).Value = 77
$ws.Cells.Item(23, 6).Value = 44
$ws.Cells.Item(24, 6).Value = 325
$ws.Cells.Item(25, 6).Value = 334
$ws.Cells.Item(26, 6).Value = 3881
$ws.Cells.Item(29, 6).Value = 40
$ws.Cells.Item(30, 6).Value = 373
$ws.Cells.Item(31, 6).Value = 957
$ws.Cells.Item(32, 6).Value = 1501
$ws.Cells.Item(34, 6).Value = 62
$ws.Cells.Item(35, 6).Value = 2775
$ws.Cells.Item(36, 6).Value = 1897
$ws.Cells.Item(37, 6).Value = 37
$ws.Cells.Item(38, 6).Value = 52
$ws.Cells.Item(39, 6).Value = 64
$ws.Cells.Item(41, 6).Value = 334
$ws.Cells.Item(42, 6).Value = 283
$ws.Cells.Item(44, 6).Value = 45
$ws.Cells.Item(45, 6).Value = 924
$ws.Cells.Item(46, 6).Value = 551
$ws.Cells.Item(47, 6).Value = 121
$ws.Cells.Item(49, 6).Value = 561
$ws.Cells.Item(50, 6).Value = 652

